# Add a new "Jira name" column (D) to the Users sheet with the two known
# values, and move the active selection to D19 (matches the author's final
# cursor position after entering the data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + data for column D
$ws.Range("D1").Value = "Jira name"
$ws.Range("D2").Value = "Joseph Steele-Perkins"
$ws.Range("D4").Value = "Developer"

# Leave the user's selection where they left off editing
$ws.Range("D19").Select()
